$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "retroalimentación de tus pares y docentes",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "retroalimentación de tus pares y docentes", 2
)
